# Auto-generated edit script: updates the cryptos price table (Sheet1)
# to reflect the latest scrape (coin swap at rows 42/43 + refreshed
# price/volume figures), per the "Updated cryptos list" GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swap (rows 42 & 43 traded ranking positions) ---
$ws.Range("B42").Value = "RenderToken"
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"

# --- Price column (D) updates ---
# Column D stores prices as plain text (e.g. "67.815.44", "0.999") so we
# force the cell to Text format before writing; otherwise Excel's COM
# layer would auto-parse values like "0.999" or "587.08" as numbers.
# The format is reset back to the default immediately afterwards so the
# cell style matches the original (unstyled) cells.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.815.44"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.493.36"
$cell.Style = "Normal"
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "587.08"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "176.95"
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.142"
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.339"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "2.948.35"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "25.69"
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "67.658.97"
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.498.10"
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "7.53"
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "350.61"
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.82"
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "4.27"
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.11"
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.619.74"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0904"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "504.44"
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "162.35"
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "18.67"
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.33"
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "4.86"
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.328"
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "144.76"
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.515"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0742"
$cell.Style = "Normal"

# --- Volume(1h) column (E) updates ---
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +3.71%  "
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +3.42%  "
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +3.34%  "
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +0.32%  "

